$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed odds values for rows 3-8 (columns G:BD)
$ws.Range("G3").Value = 2.15
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 3.4
$ws.Range("J3").Value = 2.88
$ws.Range("L3").Value = 4
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 9.5
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 19
$ws.Range("AA3").Value = 19
$ws.Range("AC3").Value = 8.5
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 9
$ws.Range("AH3").Value = 17
$ws.Range("AI3").Value = 12
$ws.Range("AJ3").Value = 34
$ws.Range("AK3").Value = 29
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 12
$ws.Range("AT3").Value = 2.63
$ws.Range("AV3").Value = 51
$ws.Range("AW3").Value = 5
$ws.Range("AX3").Value = 19
$ws.Range("AY3").Value = 29
$ws.Range("AZ3").Value = 67
$ws.Range("BB3").Value = 251
$ws.Range("G4").Value = 2.25
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 3.2
$ws.Range("M4").Value = 1.05
$ws.Range("O4").Value = 1.25
$ws.Range("Q4").Value = 1.85
$ws.Range("R4").Value = 2
$ws.Range("X4").Value = 12
$ws.Range("Y4").Value = 9.5
$ws.Range("AD4").Value = 6
$ws.Range("AG4").Value = 11
$ws.Range("AK4").Value = 23
$ws.Range("G5").Value = 1.83
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 4.5
$ws.Range("J5").Value = 2.6
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 2.4
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("X5").Value = 7.5
$ws.Range("Z5").Value = 15
$ws.Range("AA5").Value = 19
$ws.Range("AB5").Value = 41
$ws.Range("AC5").Value = 7
$ws.Range("AD5").Value = 6.5
$ws.Range("AG5").Value = 9.5
$ws.Range("AH5").Value = 21
$ws.Range("AI5").Value = 17
$ws.Range("AJ5").Value = 51
$ws.Range("AK5").Value = 41
$ws.Range("AN5").Value = 3.6
$ws.Range("AO5").Value = 10
$ws.Range("AP5").Value = 26
$ws.Range("AQ5").Value = 41
$ws.Range("AR5").Value = 67
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.38
$ws.Range("AV5").Value = 67
$ws.Range("AW5").Value = 6
$ws.Range("AX5").Value = 29
$ws.Range("AZ5").Value = 101
$ws.Range("BA5").Value = 151
$ws.Range("BB5").Value = 351
$ws.Range("G6").Value = 3.2
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 2.3
$ws.Range("K6").Value = 1.95
$ws.Range("L6").Value = 3.2
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 7
$ws.Range("O6").Value = 1.5
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.5
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("W6").Value = 7.5
$ws.Range("X6").Value = 15
$ws.Range("Y6").Value = 13
$ws.Range("AA6").Value = 34
$ws.Range("AJ6").Value = 21
$ws.Range("AO6").Value = 21
$ws.Range("BB6").Value = 251
$ws.Range("G7").Value = 2.35
$ws.Range("I7").Value = 2.9
$ws.Range("L7").Value = 3.6
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9
$ws.Range("P7").Value = 3.25
$ws.Range("AB7").Value = 34
$ws.Range("AJ7").Value = 29
$ws.Range("AY7").Value = 26
$ws.Range("G8").Value = 1.75
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 2.5
$ws.Range("L8").Value = 5.5
$ws.Range("S8").Value = 1.5
$ws.Range("T8").Value = 2.5
$ws.Range("U8").Value = 2.1
$ws.Range("V8").Value = 1.67
$ws.Range("Y8").Value = 9
$ws.Range("AB8").Value = 34
$ws.Range("AC8").Value = 7
$ws.Range("AD8").Value = 6.5
$ws.Range("AP8").Value = 23
$ws.Range("AQ8").Value = 34
$ws.Range("AS8").Value = 201
$ws.Range("AT8").Value = 2.5
$ws.Range("AU8").Value = 9

# Remove row 9 (Santos Laguna - Juarez, MEXICO - LIGA MX) entirely
$ws.Rows.Item(9).Delete()
